$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.068.41"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "3.904.18"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "464.55"
$ws.Range("E5").Value = "  +8.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.77"
$ws.Range("E6").Value = "  +4.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  +7.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000340"
$ws.Range("E11").Value = "  +8.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.95"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "4.528.68"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "3.894.88"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.98"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").Value = "  +2.86%  "
$ws.Range("D20").Value = "67.320.25"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.01"
$ws.Range("E21").Value = "  +6.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.69"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.62"
$ws.Range("E24").Value = "  +4.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "38.54"
$ws.Range("E25").Value = "  +5.20%  "
$ws.Range("E26").Value = "  +7.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.71"
$ws.Range("E27").Value = "  +5.41%  "
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "738.62"
$ws.Range("E30").Value = "  +6.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.63"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.72"
$ws.Range("E33").Value = "  -2.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.02"
$ws.Range("E34").Value = "  +6.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.158"
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.13"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "0.0₃0797"
$ws.Range("E38").Value = "  +20.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.22"
$ws.Range("E40").Value = "  +13.30%  "
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +5.58%  "
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("E46").Value = "  +5.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.40"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.90"
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.12"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.18"
$ws.Range("E51").Value = "  +0.63%  "
